# Update the "laporan kinerja" report from the Aug/Sep-2022 figures to the
# Sep/Oct-2022 figures (refreshed month labels + recomputed totals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Month/period labels -----------------------------------------------
# "Sep 2022" / "Oct 2022" look like dates to Excel's smart-entry parser, so
# force the cell format to Text ("@") first to keep them as literal text
# instead of being converted into date serial numbers.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "Sep 2022"

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "Oct 2022"

# --- Section 1: Tingkat Pengembalian Pinjaman Mitra Binaan -----------------
# These cells store numbers-with-thousands-separators as literal TEXT in the
# workbook, so force the cell format to Text ("@") before writing the value
# to stop Excel from reinterpreting the punctuated string as a real number.
$cells = @(
    @{ Ref = "C7";  Value = "1,929,693,500" },
    @{ Ref = "E7";  Value = "1,929,693,500" },
    @{ Ref = "C8";  Value = "448,029,000" },
    @{ Ref = "E8";  Value = "336,021,750" },
    @{ Ref = "C9";  Value = "1,024,203,873" },
    @{ Ref = "E9";  Value = "256,050,968" },
    @{ Ref = "C10"; Value = "1,291,025,515" },
    @{ Ref = "C11"; Value = "4,692,951,888" },
    @{ Ref = "E11"; Value = "2,521,766,218" },
    @{ Ref = "C20"; Value = "1,320,000,000" },
    @{ Ref = "C22"; Value = "1,320,000,000" },
    @{ Ref = "C28"; Value = "141,498,235" },
    @{ Ref = "C29"; Value = "5,618,783" },
    @{ Ref = "C30"; Value = "-1,662,255" },
    @{ Ref = "C31"; Value = "1,525,155,952" },
    @{ Ref = "C32"; Value = "1,997,311,245" }
)

foreach ($cell in $cells) {
    $rng = $ws.Range($cell.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cell.Value
}

# --- Narrative / summary sentences (plain text) -----------------------------
$ws.Range("A13").Value = "Kolektibilitas = 4,692,951,888 / 2521766218.25 = 53.74 %"
$ws.Range("A34").Value = "Jumlah Dana Yg Disalurkan/Jumlah Dana Tersedia= 1997311244.53 / 2582610530.00 = 66.09 %"
